$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.018.83"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -0.08%  "
$ws.Range("D3").Value = "'1.827.60"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -0.14%  "
$ws.Range("D4").Value = "'0.9979"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "'244.33"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.24%  "
$ws.Range("D6").Value = "'0.6320"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.46%  "
$ws.Range("D7").Value = "'0.9995"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").Value = "'0.07527"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -1.09%  "
$ws.Range("D9").Value = "'0.2940"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.86%  "
$ws.Range("D10").Value = "'23.07"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +1.19%  "
$ws.Range("D11").Value = "'0.07696"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.58%  "
$ws.Range("D12").Value = "'1.827.84"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -0.11%  "
$ws.Range("D13").Value = "'4.996"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +0.73%  "
$ws.Range("D14").Value = "'0.6691"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +0.55%  "
$ws.Range("D15").Value = "'83.14"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.79%  "
$ws.Range("D16").Value = "'0.000009587"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +1.81%  "
$ws.Range("D17").Value = "'6.070"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +1.27%  "
$ws.Range("D18").Value = "'29.053.30"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +1.67%  "
$ws.Range("D19").Value = "'12.59"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +2.10%  "
$ws.Range("D20").Value = "'226.72"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.55%  "
$ws.Range("D21").Value = "'0.9987"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.05%  "
$ws.Range("E22").Value = "  -1.14%  "
$ws.Range("D23").Value = "'0.9989"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.19%  "
$ws.Range("D24").Value = "'160.31"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.17%  "
$ws.Range("D25").Value = "'0.1426"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +4.24%  "
$ws.Range("D26").Value = "'8.510"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +1.02%  "
$ws.Range("D27").Value = "'17.93"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +0.48%  "
$ws.Range("D28").Value = "'1.504"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +0.76%  "
$ws.Range("D29").Value = "'4.147"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +2.39%  "
$ws.Range("D30").Value = "'4.063"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +0.73%  "
$ws.Range("D31").Value = "'0.05473"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +5.11%  "
$ws.Range("D32").Value = "'1.201"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -0.07%  "
$ws.Range("D33").Value = "'1.856"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.64%  "
$ws.Range("D34").Value = "'0.7447"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +1.78%  "
$ws.Range("D35").Value = "'1.137"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -1.55%  "
$ws.Range("D36").Value = "'2.656"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +1.70%  "
$ws.Range("D37").Value = "'1.243.84"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -2.60%  "
$ws.Range("D38").Value = "'2.752"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -0.26%  "
$ws.Range("D39").Value = "'0.01783"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.15%  "
$ws.Range("D40").Value = "'6.617"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +1.85%  "
$ws.Range("D41").Value = "'0.9022"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +1.49%  "
$ws.Range("D42").Value = "'0.9991"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -0.07%  "
$ws.Range("D43").Value = "'101.21"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.37%  "
$ws.Range("D44").Value = "'1.975.51"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +0.14%  "
$ws.Range("B45").Value = "BabyDogeCoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D45").Value = "'0.00000000124"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +2.51%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "'64.99"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +1.85%  "
$ws.Range("D47").Value = "'0.5096"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.21%  "
$ws.Range("D48").Value = "'0.4066"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +2.07%  "
$ws.Range("D49").Value = "'9.005"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +2.02%  "
$ws.Range("E50").Value = "  +0.54%  "
$ws.Range("D51").Value = "'0.05782"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.79%  "
